# Auto-generated edit script: updates Louisoix_Profits leve-profit values
# per the scheduled-runner recompute (price refresh) commit.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 11825
$ws.Range("I74").Value = 9800
$ws.Range("K74").Value = 9800
$ws.Range("M74").Value = -8864
$ws.Range("H77").Value = 11825
$ws.Range("I77").Value = 9800
$ws.Range("K77").Value = 49000
$ws.Range("M77").Value = -44320
$ws.Range("H80").Value = 1035.8889
$ws.Range("I80").Value = 1505.8
$ws.Range("J80").Value = 855.1539
$ws.Range("K80").Value = 4517.4
$ws.Range("L80").Value = 2565.4617
$ws.Range("M80").Value = -3519.4
$ws.Range("N80").Value = -4561.4617
$ws.Range("H83").Value = 1035.8889
$ws.Range("I83").Value = 1505.8
$ws.Range("J83").Value = 855.1539
$ws.Range("K83").Value = 13552.2
$ws.Range("L83").Value = 7696.3851
$ws.Range("M83").Value = -8560.199999999999
$ws.Range("N83").Value = -17680.3851
$ws.Range("H115").Value = 607
$ws.Range("I115").Value = 607
$ws.Range("K115").Value = 1821
$ws.Range("M115").Value = -254
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120
$ws.Range("H141").Value = 1522.3334
$ws.Range("I141").Value = 1522.3334
$ws.Range("K141").Value = 4567.0002
$ws.Range("M141").Value = 612.9997999999996

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 32059.25
$ws.Range("J62").Value = 32059.25
$ws.Range("L62").Value = 32059.25
$ws.Range("N62").Value = -33307.25
$ws.Range("H65").Value = 32059.25
$ws.Range("J65").Value = 32059.25
$ws.Range("L65").Value = 96177.75
$ws.Range("N65").Value = -102417.75
$ws.Range("H74").Value = 959.6
$ws.Range("I74").Value = 974
$ws.Range("K74").Value = 974
$ws.Range("M74").Value = -100
$ws.Range("H77").Value = 959.6
$ws.Range("I77").Value = 974
$ws.Range("K77").Value = 4870
$ws.Range("M77").Value = -502
$ws.Range("H109").Value = 47149.5
$ws.Range("J109").Value = 47149.5
$ws.Range("L109").Value = 47149.5
$ws.Range("N109").Value = -49923.5
$ws.Range("H114").Value = 93166.336
$ws.Range("J114").Value = 93166.336
$ws.Range("L114").Value = 93166.336
$ws.Range("N114").Value = -101844.336

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18765.428
$ws.Range("J82").Value = 70000
$ws.Range("L82").Value = 70000
$ws.Range("N82").Value = -70766
$ws.Range("H85").Value = 18765.428
$ws.Range("J85").Value = 70000
$ws.Range("L85").Value = 70000
$ws.Range("N85").Value = -72652
$ws.Range("H112").Value = 127632
$ws.Range("J112").Value = 127632
$ws.Range("L112").Value = 127632
$ws.Range("N112").Value = -130586
$ws.Range("H134").Value = 2283.7292
$ws.Range("I134").Value = 1761.5952
$ws.Range("J134").Value = 5938.6665
$ws.Range("K134").Value = 5284.7856
$ws.Range("L134").Value = 17815.9995
$ws.Range("M134").Value = -2749.7856
$ws.Range("N134").Value = -22885.9995

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 999
$ws.Range("J4").Value = 999
$ws.Range("L4").Value = 999
$ws.Range("N4").Value = -1223
$ws.Range("H58").Value = 79641.69500000001
$ws.Range("I58").Value = 113667.336
$ws.Range("J58").Value = 3084
$ws.Range("K58").Value = 113667.336
$ws.Range("L58").Value = 3084
$ws.Range("M58").Value = -113464.336
$ws.Range("N58").Value = -3490
$ws.Range("H69").Value = 77499.5
$ws.Range("I69").Value = 50000
$ws.Range("J69").Value = 104999
$ws.Range("K69").Value = 50000
$ws.Range("L69").Value = 104999
$ws.Range("M69").Value = -49251
$ws.Range("N69").Value = -106497
$ws.Range("H72").Value = 77499.5
$ws.Range("I72").Value = 50000
$ws.Range("J72").Value = 104999
$ws.Range("K72").Value = 150000
$ws.Range("L72").Value = 314997
$ws.Range("M72").Value = -146256
$ws.Range("N72").Value = -322485
$ws.Range("H107").Value = 2293.9143
$ws.Range("I107").Value = 725.4375
$ws.Range("K107").Value = 725.4375
$ws.Range("M107").Value = 1194.5625
$ws.Range("H134").Value = 81898.734
$ws.Range("J134").Value = 10504.667
$ws.Range("L134").Value = 31514.001
$ws.Range("N134").Value = -36584.001
$ws.Range("H136").Value = 79641.69500000001
$ws.Range("I136").Value = 113667.336
$ws.Range("J136").Value = 3084
$ws.Range("K136").Value = 341002.008
$ws.Range("L136").Value = 9252
$ws.Range("M136").Value = -338452.008
$ws.Range("N136").Value = -14352

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2326.3333
$ws.Range("J68").Value = 4999
$ws.Range("L68").Value = 14997
$ws.Range("N68").Value = -16619
$ws.Range("H71").Value = 2326.3333
$ws.Range("J71").Value = 4999
$ws.Range("L71").Value = 44991
$ws.Range("N71").Value = -53103
$ws.Range("H81").Value = 8533.333000000001
$ws.Range("J81").Value = 12500
$ws.Range("L81").Value = 37500
$ws.Range("N81").Value = -39746
$ws.Range("H84").Value = 8533.333000000001
$ws.Range("J84").Value = 12500
$ws.Range("L84").Value = 112500
$ws.Range("N84").Value = -123732
$ws.Range("H86").Value = 439.125
$ws.Range("J86").Value = 605
$ws.Range("L86").Value = 1815
$ws.Range("N86").Value = -4187
$ws.Range("H87").Value = 4912.25
$ws.Range("I87").Value = 4912.25
$ws.Range("K87").Value = 14736.75
$ws.Range("M87").Value = -13488.75
$ws.Range("H88").Value = 7999
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 439.125
$ws.Range("J89").Value = 605
$ws.Range("L89").Value = 5445
$ws.Range("N89").Value = -17301
$ws.Range("H90").Value = 4912.25
$ws.Range("I90").Value = 4912.25
$ws.Range("K90").Value = 44210.25
$ws.Range("M90").Value = -37970.25
$ws.Range("H91").Value = 7999
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 79996.38
$ws.Range("I113").Value = 79551
$ws.Range("K113").Value = 79551
$ws.Range("M113").Value = -77381
$ws.Range("H118").Value = 42499
$ws.Range("J118").Value = 42499
$ws.Range("L118").Value = 42499
$ws.Range("N118").Value = -45813
$ws.Range("H132").Value = 55091.316
$ws.Range("I132").Value = 61290.35
$ws.Range("J132").Value = 2399.5
$ws.Range("K132").Value = 183871.05
$ws.Range("L132").Value = 7198.5
$ws.Range("M132").Value = -181341.05
$ws.Range("N132").Value = -12258.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 39648.75
$ws.Range("J42").Value = 39699.5
$ws.Range("L42").Value = 39699.5
$ws.Range("N42").Value = -40825.5
$ws.Range("H45").Value = 50000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 4644.8
$ws.Range("J46").Value = 2928.0378
$ws.Range("L46").Value = 2928.0378
$ws.Range("N46").Value = -3304.0378
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H49").Value = 39648.75
$ws.Range("J49").Value = 39699.5
$ws.Range("L49").Value = 39699.5
$ws.Range("N49").Value = -39993.5
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 339323
$ws.Range("J62").Value = 339323
$ws.Range("L62").Value = 339323
$ws.Range("N62").Value = -340571
$ws.Range("H65").Value = 339323
$ws.Range("J65").Value = 339323
$ws.Range("L65").Value = 1696615
$ws.Range("N65").Value = -1702855

